$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row layout: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric metric columns (E..T).
$rows = @(
    @("M1", "Matn1", "Itga1", "ECs", 1, 0.3333333333333333, 0.306399, 0.919197, 0.3017954127735523, 0.3017954127735522, 2, 1, 55.9484575, 111.896915, 0.4236195045309116, 0.3537406883573876, 17.1425514295425, 102.855308577255, 0.1278464232288342, 0.1067573170576183),
    @("M1", "Matn1", "Itga1", "FAPs", 1, 0.3333333333333333, 0.306399, 0.919197, 0.3017954127735523, 0.3017954127735522, 3, 1, 2.360452, 7.081356, 0.01787240527067255, 0.02238635217015337, 0.7232401323480001, 6.509161191132, 0.005393809925918834, 0.006756098393685541),
    @("M1", "Matn1", "Itga1", "M1", 1, 0.3333333333333333, 0.306399, 0.919197, 0.3017954127735523, 0.3017954127735522, 3, 1, 2.631462, 7.894386, 0.01992438537973851, 0.02495659096409337, 0.8062773253380001, 7.256495928042, 0.006013088109937515, 0.007531784671429262),
    @("M1", "Matn1", "Itga1", "M2", 1, 0.3333333333333333, 0.306399, 0.919197, 0.3017954127735523, 0.3017954127735522, 3, 1, 0.4270753333333333, 1.281226, 0.003233644843632026, 0.004050350871437182, 0.130855455058, 1.177699095522, 0.0009758991803469961, 0.001222377313123101),
    @("M1", "Matn1", "Itga1", "Neutro", 1, 0.3333333333333333, 0.306399, 0.919197, 0.3017954127735523, 0.3017954127735522, 3, 1, 46.760844, 140.282532, 0.3540545432682795, 0.4434763856912163, 14.327475840756, 128.947282566804, 0.1068520370300019, 0.1338391388750037),
    @("M1", "Matn1", "Itga1", "sCs", 1, 0.3333333333333333, 0.306399, 0.919197, 0.3017954127735523, 0.3017954127735522, 2, 1, 23.9441395, 47.888279, 0.1812955167067658, 0.1513896319457121, 7.3364603986605, 44.018762391963, 0.05471415529851283, 0.04568869646269232),
    @("M2", "Matn1", "Itga1", "ECs", 3, 1, 0.7088550000000001, 2.126565, 0.6982045872264478, 0.6982045872264477, 2, 1, 55.9484575, 111.896915, 0.4236195045309116, 0.3537406883573876, 39.65934384116251, 237.956063046975, 0.2957730813020775, 0.2469833712997693),
    @("M2", "Matn1", "Itga1", "FAPs", 3, 1, 0.7088550000000001, 2.126565, 0.6982045872264478, 0.6982045872264477, 3, 1, 2.360452, 7.081356, 0.01787240527067255, 0.02238635217015337, 1.67321820246, 15.05896382214, 0.01247859534475372, 0.01563025377646782),
    @("M2", "Matn1", "Itga1", "M1", 3, 1, 0.7088550000000001, 2.126565, 0.6982045872264478, 0.6982045872264477, 3, 1, 2.631462, 7.894386, 0.01992438537973851, 0.02495659096409337, 1.86532499601, 16.78792496409, 0.013911297269801, 0.01742480629266411),
    @("M2", "Matn1", "Itga1", "M2", 3, 1, 0.7088550000000001, 2.126565, 0.6982045872264478, 0.6982045872264477, 3, 1, 0.4270753333333333, 1.281226, 0.003233644843632026, 0.004050350871437182, 0.30273448541, 2.72461036869, 0.00225774566328503, 0.002827973558314081),
    @("M2", "Matn1", "Itga1", "Neutro", 3, 1, 0.7088550000000001, 2.126565, 0.6982045872264478, 0.6982045872264477, 3, 1, 46.760844, 140.282532, 0.3540545432682795, 0.4434763856912163, 33.14665807362, 298.31992266258, 0.2472025062382776, 0.3096372468162126),
    @("M2", "Matn1", "Itga1", "sCs", 3, 1, 0.7088550000000001, 2.126565, 0.6982045872264478, 0.6982045872264477, 2, 1, 23.9441395, 47.888279, 0.1812955167067658, 0.1513896319457121, 16.9729230052725, 101.837538031635, 0.126581361408253, 0.1057009354830197)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}